# Auto-generated edit script applying numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 61666.332
$ws.Range("J3").Value = 61666.332
$ws.Range("L3").Value = 61666.332
$ws.Range("N3").Value = -61894.332
$ws.Range("H33").Value = 824.5909
$ws.Range("I33").Value = 367.63635
$ws.Range("K33").Value = 367.63635
$ws.Range("M33").Value = -138.63635
$ws.Range("H40").Value = 5026.727
$ws.Range("I40").Value = 5220.7144
$ws.Range("K40").Value = 5220.7144
$ws.Range("M40").Value = -5045.7144
$ws.Range("H41").Value = 1368.3572
$ws.Range("J41").Value = 1619.7273
$ws.Range("L41").Value = 1619.7273
$ws.Range("N41").Value = -2499.7273
$ws.Range("H86").Value = 2032.25
$ws.Range("I86").Value = 1959.25
$ws.Range("J86").Value = 2105.25
$ws.Range("K86").Value = 1959.25
$ws.Range("L86").Value = 2105.25
$ws.Range("M86").Value = -836.25
$ws.Range("N86").Value = -4351.25
$ws.Range("H89").Value = 2032.25
$ws.Range("I89").Value = 1959.25
$ws.Range("J89").Value = 2105.25
$ws.Range("K89").Value = 9796.25
$ws.Range("L89").Value = 10526.25
$ws.Range("M89").Value = -4180.25
$ws.Range("N89").Value = -21758.25
$ws.Range("H102").Value = 61666.332
$ws.Range("J102").Value = 61666.332
$ws.Range("L102").Value = 61666.332
$ws.Range("N102").Value = -68156.33199999999
$ws.Range("H132").Value = 3191529.8
$ws.Range("I132").Value = 3264021
$ws.Range("K132").Value = 9792063
$ws.Range("M132").Value = -9789533

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14858.343
$ws.Range("I32").Value = 15526.942
$ws.Range("K32").Value = 15526.942
$ws.Range("M32").Value = -15239.942
$ws.Range("H45").Value = 2294.3809
$ws.Range("J45").Value = 4362.2856
$ws.Range("L45").Value = 4362.2856
$ws.Range("N45").Value = -5116.2856
$ws.Range("H74").Value = 120534.78
$ws.Range("I74").Value = 143519.05
$ws.Range("J74").Value = 13274.889
$ws.Range("K74").Value = 143519.05
$ws.Range("L74").Value = 13274.889
$ws.Range("M74").Value = -142645.05
$ws.Range("N74").Value = -15022.889
$ws.Range("H77").Value = 120534.78
$ws.Range("I77").Value = 143519.05
$ws.Range("J77").Value = 13274.889
$ws.Range("K77").Value = 717595.25
$ws.Range("L77").Value = 66374.44499999999
$ws.Range("M77").Value = -713227.25
$ws.Range("N77").Value = -75110.44499999999
$ws.Range("H132").Value = 1107.5745
$ws.Range("I132").Value = 946.5952
$ws.Range("K132").Value = 2839.7856
$ws.Range("M132").Value = -309.7856000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14274.292
$ws.Range("I20").Value = 18587.777
$ws.Range("J20").Value = 1333.8334
$ws.Range("K20").Value = 18587.777
$ws.Range("L20").Value = 1333.8334
$ws.Range("M20").Value = -18340.777
$ws.Range("N20").Value = -1827.8334
$ws.Range("H86").Value = 1349
$ws.Range("I86").Value = 1284.909
$ws.Range("K86").Value = 1284.909
$ws.Range("M86").Value = -161.9090000000001
$ws.Range("H89").Value = 1349
$ws.Range("I89").Value = 1284.909
$ws.Range("K89").Value = 6424.545
$ws.Range("M89").Value = -808.5450000000001
$ws.Range("H107").Value = 724.1111
$ws.Range("I107").Value = 689.75
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 689.75
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1230.25
$ws.Range("N107").Value = -4839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 674.1111
$ws.Range("I22").Value = 326.33334
$ws.Range("J22").Value = 848
$ws.Range("K22").Value = 326.33334
$ws.Range("L22").Value = 848
$ws.Range("M22").Value = 23.66665999999998
$ws.Range("N22").Value = -1548
$ws.Range("H58").Value = 8625.977999999999
$ws.Range("I58").Value = 1163.775
$ws.Range("K58").Value = 1163.775
$ws.Range("M58").Value = -960.7750000000001
$ws.Range("H86").Value = 14391.333
$ws.Range("I86").Value = 17853.908
$ws.Range("J86").Value = 10582.5
$ws.Range("K86").Value = 17853.908
$ws.Range("L86").Value = 10582.5
$ws.Range("M86").Value = -16730.908
$ws.Range("N86").Value = -12828.5
$ws.Range("H89").Value = 14391.333
$ws.Range("I89").Value = 17853.908
$ws.Range("J89").Value = 10582.5
$ws.Range("K89").Value = 89269.53999999999
$ws.Range("L89").Value = 52912.5
$ws.Range("M89").Value = -83653.53999999999
$ws.Range("N89").Value = -64144.5
$ws.Range("H94").Value = 1649.6086
$ws.Range("J94").Value = 1806.5714
$ws.Range("L94").Value = 1806.5714
$ws.Range("N94").Value = -2708.5714
$ws.Range("H99").Value = 5732.1055
$ws.Range("I99").Value = 4026.7334
$ws.Range("K99").Value = 4026.7334
$ws.Range("M99").Value = -2528.7334
$ws.Range("H107").Value = 1999.4
$ws.Range("I107").Value = 1999
$ws.Range("K107").Value = 1999
$ws.Range("M107").Value = -79
$ws.Range("H122").Value = 1492.5312
$ws.Range("J122").Value = 1708.7142
$ws.Range("L122").Value = 5126.142599999999
$ws.Range("N122").Value = -10026.1426
$ws.Range("H126").Value = 5732.1055
$ws.Range("I126").Value = 4026.7334
$ws.Range("K126").Value = 12080.2002
$ws.Range("M126").Value = -9610.200199999999
$ws.Range("H136").Value = 8625.977999999999
$ws.Range("I136").Value = 1163.775
$ws.Range("K136").Value = 3491.325
$ws.Range("M136").Value = -941.3250000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1307.6
$ws.Range("I109").Value = 180
$ws.Range("K109").Value = 540
$ws.Range("M109").Value = 500
$ws.Range("H122").Value = 786.3333
$ws.Range("J122").Value = 878.86957
$ws.Range("L122").Value = 7909.826129999999
$ws.Range("N122").Value = -12809.82613

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 183.6
$ws.Range("J2").Value = 194.66667
$ws.Range("L2").Value = 194.66667
$ws.Range("N2").Value = -420.66667
$ws.Range("H97").Value = 1640.9259
$ws.Range("I97").Value = 1156.1111
$ws.Range("J97").Value = 2610.5557
$ws.Range("K97").Value = 1156.1111
$ws.Range("L97").Value = 2610.5557
$ws.Range("M97").Value = -660.1111000000001
$ws.Range("N97").Value = -3602.5557
$ws.Range("H107").Value = 416.33334
$ws.Range("I107").Value = 99.666664
$ws.Range("K107").Value = 99.666664
$ws.Range("M107").Value = 1820.333336
$ws.Range("H122").Value = 2533.7188
$ws.Range("I122").Value = 2180.2
$ws.Range("K122").Value = 6540.599999999999
$ws.Range("M122").Value = -4090.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3570.2856
$ws.Range("I7").Value = 3848.125
$ws.Range("K7").Value = 3848.125
$ws.Range("M7").Value = -3736.125
$ws.Range("H122").Value = 6507.5454
$ws.Range("I122").Value = 4900
$ws.Range("K122").Value = 14700
$ws.Range("M122").Value = -12250
$ws.Range("H126").Value = 3570.2856
$ws.Range("I126").Value = 3848.125
$ws.Range("K126").Value = 11544.375
$ws.Range("M126").Value = -9074.375
$ws.Range("H136").Value = 3165.2
$ws.Range("I136").Value = 2930
$ws.Range("J136").Value = 3870.8
$ws.Range("K136").Value = 8790
$ws.Range("L136").Value = 11612.4
$ws.Range("M136").Value = -6240
$ws.Range("N136").Value = -16712.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 50267.703
$ws.Range("I122").Value = 53409.12
$ws.Range("K122").Value = 160227.36
$ws.Range("M122").Value = -157777.36
$ws.Range("H136").Value = 26347.545
$ws.Range("I136").Value = 27197.477
$ws.Range("K136").Value = 81592.431
$ws.Range("M136").Value = -79042.431

